$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H55").Value = 671
$ws_ALC.Range("J55").Value = 541.2
$ws_ALC.Range("L55").Value = 541.2
$ws_ALC.Range("N55").Value = -969.2

$ws_ALC.Range("H82").Value = 4303.5713
$ws_ALC.Range("I82").Value = 2530
$ws_ALC.Range("J82").Value = 6668.3335
$ws_ALC.Range("K82").Value = 7590
$ws_ALC.Range("L82").Value = 20005.0005
$ws_ALC.Range("M82").Value = -7184
$ws_ALC.Range("N82").Value = -20817.0005

$ws_ALC.Range("H85").Value = 4303.5713
$ws_ALC.Range("I85").Value = 2530
$ws_ALC.Range("J85").Value = 6668.3335
$ws_ALC.Range("K85").Value = 7590
$ws_ALC.Range("L85").Value = 20005.0005
$ws_ALC.Range("M85").Value = -6186
$ws_ALC.Range("N85").Value = -22813.0005

$ws_ALC.Range("H116").Value = 5424150.5
$ws_ALC.Range("I116").Value = 10842172
$ws_ALC.Range("J116").Value = 6129
$ws_ALC.Range("K116").Value = 10842172
$ws_ALC.Range("L116").Value = 6129
$ws_ALC.Range("M116").Value = -10838730
$ws_ALC.Range("N116").Value = -13013

$ws_ALC.Range("H132").Value = 11111.47
$ws_ALC.Range("I132").Value = 2030.7778
$ws_ALC.Range("K132").Value = 6092.3334
$ws_ALC.Range("M132").Value = -3562.3334

$ws_ALC.Range("H136").Value = 122390
$ws_ALC.Range("J136").Value = 122390
$ws_ALC.Range("L136").Value = 122390
$ws_ALC.Range("N136").Value = -132590

$ws_ALC.Range("H137").Value = 10421286
$ws_ALC.Range("I137").Value = 1022.125
$ws_ALC.Range("K137").Value = 3066.375
$ws_ALC.Range("M137").Value = -516.375

$ws_ALC.Range("H139").Value = 129834.75
$ws_ALC.Range("J139").Value = 129834.75
$ws_ALC.Range("L139").Value = 129834.75
$ws_ALC.Range("N139").Value = -140114.75

$ws_ALC.Range("H140").Value = 61759.777
$ws_ALC.Range("J140").Value = 60641.125
$ws_ALC.Range("L140").Value = 60641.125
$ws_ALC.Range("N140").Value = -71001.125

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H22").Value = 477.44446
$ws_BSM.Range("I22").Value = 477.44446
$ws_BSM.Range("K22").Value = 477.44446
$ws_BSM.Range("M22").Value = -304.44446

$ws_BSM.Range("H99").Value = 802289.1
$ws_BSM.Range("I99").Value = 1158312.4
$ws_BSM.Range("K99").Value = 1158312.4
$ws_BSM.Range("M99").Value = -1156814.4

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 21744734
$ws_CRP.Range("I31").Value = 43480444
$ws_CRP.Range("J31").Value = 9023.652
$ws_CRP.Range("K31").Value = 43480444
$ws_CRP.Range("L31").Value = 9023.652
$ws_CRP.Range("M31").Value = -43480149
$ws_CRP.Range("N31").Value = -9613.652

$ws_CRP.Range("H34").Value = 21744734
$ws_CRP.Range("I34").Value = 43480444
$ws_CRP.Range("J34").Value = 9023.652
$ws_CRP.Range("K34").Value = 43480444
$ws_CRP.Range("L34").Value = 9023.652
$ws_CRP.Range("M34").Value = -43480242
$ws_CRP.Range("N34").Value = -9427.652

$ws_CRP.Range("H141").Value = 96432.08
$ws_CRP.Range("J141").Value = 96432.08
$ws_CRP.Range("L141").Value = 96432.08
$ws_CRP.Range("N141").Value = -106792.08

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H131").Value = 8199301
$ws_CUL.Range("I131").Value = 13890925
$ws_CUL.Range("J131").Value = 6805434
$ws_CUL.Range("K131").Value = 41672775
$ws_CUL.Range("L131").Value = 20416302
$ws_CUL.Range("M131").Value = -41667735
$ws_CUL.Range("N131").Value = -20426382

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H97").Value = 275.1
$ws_GSM.Range("I97").Value = 413.33334
$ws_GSM.Range("J97").Value = 215.85715
$ws_GSM.Range("K97").Value = 413.33334
$ws_GSM.Range("L97").Value = 215.85715
$ws_GSM.Range("M97").Value = 82.66665999999998
$ws_GSM.Range("N97").Value = -1207.85715

$ws_GSM.Range("H122").Value = 5944
$ws_GSM.Range("I122").Value = 2970
$ws_GSM.Range("K122").Value = 8910
$ws_GSM.Range("M122").Value = -6460

$ws_GSM.Range("H140").Value = 78780
$ws_GSM.Range("J140").Value = 78780
$ws_GSM.Range("L140").Value = 78780
$ws_GSM.Range("N140").Value = -89140

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 5443
$ws_LTW.Range("I7").Value = 1971
$ws_LTW.Range("J7").Value = 7427
$ws_LTW.Range("K7").Value = 1971
$ws_LTW.Range("L7").Value = 7427
$ws_LTW.Range("M7").Value = -1859
$ws_LTW.Range("N7").Value = -7651

$ws_LTW.Range("H46").Value = 6843.615
$ws_LTW.Range("J46").Value = 7057.36
$ws_LTW.Range("L46").Value = 7057.36
$ws_LTW.Range("N46").Value = -7433.36

$ws_LTW.Range("H61").Value = 1904
$ws_LTW.Range("I61").Value = 1904
$ws_LTW.Range("K61").Value = 1904
$ws_LTW.Range("M61").Value = -1702

$ws_LTW.Range("H82").Value = 3472897.5
$ws_LTW.Range("I82").Value = 3906822.2
$ws_LTW.Range("K82").Value = 3906822.2
$ws_LTW.Range("M82").Value = -3906461.2

$ws_LTW.Range("H85").Value = 3472897.5
$ws_LTW.Range("I85").Value = 3906822.2
$ws_LTW.Range("K85").Value = 3906822.2
$ws_LTW.Range("M85").Value = -3905574.2

$ws_LTW.Range("H100").Value = 4831.8335
$ws_LTW.Range("I100").Value = 4798.5
$ws_LTW.Range("K100").Value = 4798.5
$ws_LTW.Range("M100").Value = -4257.5

$ws_LTW.Range("H113").Value = 1904
$ws_LTW.Range("I113").Value = 1904
$ws_LTW.Range("K113").Value = 1904
$ws_LTW.Range("M113").Value = 266

$ws_LTW.Range("H122").Value = 5501601.5
$ws_LTW.Range("I122").Value = 3709.75
$ws_LTW.Range("K122").Value = 11129.25
$ws_LTW.Range("M122").Value = -8679.25

$ws_LTW.Range("H126").Value = 5443
$ws_LTW.Range("I126").Value = 1971
$ws_LTW.Range("J126").Value = 7427
$ws_LTW.Range("K126").Value = 5913
$ws_LTW.Range("L126").Value = 22281
$ws_LTW.Range("M126").Value = -3443
$ws_LTW.Range("N126").Value = -27221

$ws_LTW.Range("H132").Value = 2665.3171
$ws_LTW.Range("I132").Value = 1805
$ws_LTW.Range("J132").Value = 4518.3076
$ws_LTW.Range("K132").Value = 5415
$ws_LTW.Range("L132").Value = 13554.9228
$ws_LTW.Range("M132").Value = -2885
$ws_LTW.Range("N132").Value = -18614.9228

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H122").Value = 4776.811
$ws_WVR.Range("I122").Value = 4075.52
$ws_WVR.Range("J122").Value = 6237.8335
$ws_WVR.Range("K122").Value = 12226.56
$ws_WVR.Range("L122").Value = 18713.5005
$ws_WVR.Range("M122").Value = -9776.559999999999
$ws_WVR.Range("N122").Value = -23613.5005

$ws_WVR.Range("H126").Value = 3128.3845
$ws_WVR.Range("I126").Value = 3490.2856
$ws_WVR.Range("J126").Value = 2706.1667
$ws_WVR.Range("K126").Value = 10470.8568
$ws_WVR.Range("L126").Value = 8118.500100000001
$ws_WVR.Range("M126").Value = -8000.856800000001
$ws_WVR.Range("N126").Value = -13058.5001

$ws_WVR.Range("H136").Value = 8428.743
$ws_WVR.Range("I136").Value = 4092.9375
$ws_WVR.Range("J136").Value = 11444.956
$ws_WVR.Range("K136").Value = 12278.8125
$ws_WVR.Range("L136").Value = 34334.868
$ws_WVR.Range("M136").Value = -9728.8125
$ws_WVR.Range("N136").Value = -39434.868
